# Add 2022-Q4 sheet + update 总计 summary sheet.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q4" right before the existing "2022-Q3"
#    sheet (mirrors Excel's Worksheets.Add(Before:=...) behaviour, which is
#    what shifts every sheet after it one slot to the right).
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$newSheet = $wb.Worksheets.Add($q3)
$newSheet.Name = "2022-Q4"

# Copy the header / index-column formatting (bold, border, centered) from the
# neighbouring "2022-Q3" sheet so the new sheet re-uses the same style (s="2")
# instead of minting a fresh one.
$q3.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats
$q3.Range("A2").Copy()
$newSheet.Range("A2:A12").PasteSpecial(-4122)  # xlPasteFormats

# Header row.
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Columns B, D, E, F, G hold numeric-looking figures that are stored as plain
# text in this workbook (e.g. fund codes with leading zeros). Force the
# number format to Text first so Excel doesn't silently convert them.
$newSheet.Range("B2:B12").NumberFormat = "@"
$newSheet.Range("D2:G12").NumberFormat = "@"

$q4rows = @(
  ,@("0","010963","信澳周期动力混合A","24.33","92.85","5.27","1.2822","2")
  ,@("1","015455","信澳周期动力混合C","5.15","92.85","5.27","0.2714","2")
  ,@("2","002317","招商睿逸稳健配置混合","4.84","49.57","2.98","0.1442","7")
  ,@("3","217002","招商安泰平衡混合","5.50","49.54","2.55","0.1402","7")
  ,@("4","009048","浦银安盛科技创新优选三年封闭运作灵活配置混合","2.39","35.76","1.15","0.0275","7")
  ,@("5","002292","诺安益鑫灵活配置混合A","0.37","69.58","4.22","0.0156","6")
  ,@("6","001231","银华泰利灵活配置混合A","0.87","22.82","0.78","0.0068","6")
  ,@("7","003308","中信建投睿利灵活配置混合A","0.07","71.26","3.42","0.0024","5")
  ,@("8","004635","中信建投睿利灵活配置混合C","0.04","71.26","3.42","0.0014","5")
  ,@("9","014550","诺安益鑫灵活配置混合C","0.02","69.58","4.22","0.0008","6")
  ,@("10","002328","银华泰利灵活配置混合C","0.03","22.82","0.78","0.0002","6")
)

for ($i = 0; $i -lt $q4rows.Length; $i++) {
  $row = $q4rows[$i]
  $r = $i + 2
  $newSheet.Cells.Item($r, 1).Value = [int]$row[0]     # A: index (number)
  $newSheet.Cells.Item($r, 2).Value = $row[1]           # B: fund code (text)
  $newSheet.Cells.Item($r, 3).Value = $row[2]           # C: fund name (text)
  $newSheet.Cells.Item($r, 4).Value = $row[3]           # D: fund size (text)
  $newSheet.Cells.Item($r, 5).Value = $row[4]           # E: stock position (text)
  $newSheet.Cells.Item($r, 6).Value = $row[5]           # F: position pct (text)
  $newSheet.Cells.Item($r, 7).Value = $row[6]           # G: market value (text)
  $newSheet.Cells.Item($r, 8).Value = [int]$row[7]      # H: rank (number)
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" summary sheet: a new row for 2022-Q4 is inserted right
#    after the header, pushing every existing quarter down by one row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Give the new last row (A9) the same index-column style as the rest of
# column A before writing into it.
$summary.Range("A2").Copy()
$summary.Range("A9").PasteSpecial(-4122)  # xlPasteFormats

$summaryRows = @(
  ,@("2022-Q4", 11, 1.89)
  ,@("2022-Q3", 68, 9.81)
  ,@("2022-Q2", 110, 31.57)
  ,@("2022-Q1", 62, 20.73)
  ,@("2021-Q4", 56, 18.52)
  ,@("2021-Q3", 21, 12.55)
  ,@("2021-Q2", 6, 5.11)
  ,@("2021-Q1", 1, 0.02)
)

for ($i = 0; $i -lt $summaryRows.Length; $i++) {
  $row = $summaryRows[$i]
  $r = $i + 2
  $summary.Cells.Item($r, 1).Value = $i        # A: sequential index (number)
  $summary.Cells.Item($r, 2).Value = $row[0]   # B: quarter label (text)
  $summary.Cells.Item($r, 3).Value = $row[1]   # C: holding count (number)
  $summary.Cells.Item($r, 4).Value = $row[2]   # D: holding value (number)
}
